# Append/update run: refresh the "取得日時" (fetched-at) timestamp column
# on the "ランサーズ" sheet for every existing data row (rows 2-19) to the
# new run time 2025-09-10 18:31:20.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

$newTimestamp = "2025-09-10 18:31:20"

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
if ($lastRow -lt 2) {
    $lastRow = 19
}

for ($row = 2; $row -le $lastRow; $row++) {
    $ws.Cells.Item($row, 1).Value = $newTimestamp
}
